$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 874
$ws.Range("I18").Value = 874
$ws.Range("K18").Value = 874
$ws.Range("M18").Value = -590
$ws.Range("H39").Value = 269.0909
$ws.Range("I39").Value = 246
$ws.Range("K39").Value = 738
$ws.Range("M39").Value = -442
$ws.Range("H55").Value = 350.66666
$ws.Range("I55").Value = 450
$ws.Range("J55").Value = 301
$ws.Range("K55").Value = 450
$ws.Range("L55").Value = 301
$ws.Range("M55").Value = -236
$ws.Range("N55").Value = -729
$ws.Range("H70").Value = 3135.818
$ws.Range("J70").Value = 3350.4
$ws.Range("L70").Value = 10051.2
$ws.Range("N70").Value = -10591.2
$ws.Range("H73").Value = 3135.818
$ws.Range("J73").Value = 3350.4
$ws.Range("L73").Value = 10051.2
$ws.Range("N73").Value = -11923.2
$ws.Range("H103").Value = 1799
$ws.Range("J103").Value = 1799
$ws.Range("L103").Value = 5397
$ws.Range("N103").Value = -6569
$ws.Range("H104").Value = 195.75
$ws.Range("I104").Value = 195.75
$ws.Range("K104").Value = 587.25
$ws.Range("M104").Value = 1159.75
$ws.Range("H135").Value = 914.5
$ws.Range("I135").Value = 914.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8230.5
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -5695.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 930.8333
$ws.Range("I2").Value = 841
$ws.Range("J2").Value = 1056.6
$ws.Range("K2").Value = 841
$ws.Range("L2").Value = 1056.6
$ws.Range("M2").Value = -728
$ws.Range("N2").Value = -1282.6
$ws.Range("H116").Value = 930.8333
$ws.Range("I116").Value = 841
$ws.Range("J116").Value = 1056.6
$ws.Range("K116").Value = 841
$ws.Range("L116").Value = 1056.6
$ws.Range("M116").Value = 1453
$ws.Range("N116").Value = -5644.6
$ws.Range("H132").Value = 1613.65
$ws.Range("I132").Value = 1571.6
$ws.Range("J132").Value = 1739.8
$ws.Range("K132").Value = 4714.799999999999
$ws.Range("L132").Value = 5219.4
$ws.Range("M132").Value = -2184.799999999999
$ws.Range("N132").Value = -10279.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 930.8333
$ws.Range("I3").Value = 841
$ws.Range("J3").Value = 1056.6
$ws.Range("K3").Value = 841
$ws.Range("L3").Value = 1056.6
$ws.Range("M3").Value = -727
$ws.Range("N3").Value = -1284.6
$ws.Range("H20").Value = 4678.75
$ws.Range("I20").Value = 3488.4
$ws.Range("K20").Value = 3488.4
$ws.Range("M20").Value = -3241.4
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = 0
$ws.Range("H88").Value = 13998.5
$ws.Range("J88").Value = 13998.5
$ws.Range("L88").Value = 13998.5
$ws.Range("N88").Value = -14810.5
$ws.Range("H91").Value = 13998.5
$ws.Range("J91").Value = 13998.5
$ws.Range("L91").Value = 13998.5
$ws.Range("N91").Value = -16806.5
$ws.Range("H94").Value = 3444.75
$ws.Range("I94").Value = 3444.75
$ws.Range("K94").Value = 3444.75
$ws.Range("M94").Value = -2993.75
$ws.Range("H99").Value = 7900
$ws.Range("I99").Value = 7900
$ws.Range("K99").Value = 7900
$ws.Range("M99").Value = -6402
$ws.Range("H107").Value = 1154.5
$ws.Range("I107").Value = 1109
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1109
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 811
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1883.75
$ws.Range("I31").Value = 1429.9756
$ws.Range("J31").Value = 3124.0667
$ws.Range("K31").Value = 1429.9756
$ws.Range("L31").Value = 3124.0667
$ws.Range("M31").Value = -1134.9756
$ws.Range("N31").Value = -3714.0667
$ws.Range("H34").Value = 1883.75
$ws.Range("I34").Value = 1429.9756
$ws.Range("J34").Value = 3124.0667
$ws.Range("K34").Value = 1429.9756
$ws.Range("L34").Value = 3124.0667
$ws.Range("M34").Value = -1227.9756
$ws.Range("N34").Value = -3528.0667
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H58").Value = 5576.5713
$ws.Range("I58").Value = 4856.0835
$ws.Range("K58").Value = 4856.0835
$ws.Range("M58").Value = -4653.0835
$ws.Range("H105").Value = 1439.8
$ws.Range("I105").Value = 1233
$ws.Range("K105").Value = 1233
$ws.Range("M105").Value = 514
$ws.Range("H107").Value = 1305.3334
$ws.Range("I107").Value = 1111
$ws.Range("K107").Value = 1111
$ws.Range("M107").Value = 809
$ws.Range("H132").Value = 1984.2142
$ws.Range("I132").Value = 1731.6666
$ws.Range("K132").Value = 5194.9998
$ws.Range("M132").Value = -2664.9998
$ws.Range("H136").Value = 5576.5713
$ws.Range("I136").Value = 4856.0835
$ws.Range("K136").Value = 14568.2505
$ws.Range("M136").Value = -12018.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 6701.8
$ws.Range("I6").Value = 502.33334
$ws.Range("K6").Value = 1507.00002
$ws.Range("M6").Value = -1394.00002
$ws.Range("H36").Value = 312
$ws.Range("I36").Value = 316
$ws.Range("K36").Value = 948
$ws.Range("M36").Value = -779
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("N41").Value = 0
$ws.Range("H68").Value = 1348.8
$ws.Range("J68").Value = 1500
$ws.Range("L68").Value = 4500
$ws.Range("N68").Value = -6122
$ws.Range("H71").Value = 1348.8
$ws.Range("J71").Value = 1500
$ws.Range("L71").Value = 13500
$ws.Range("N71").Value = -21612
$ws.Range("H80").Value = 12168.454
$ws.Range("J80").Value = 12098.8
$ws.Range("L80").Value = 36296.39999999999
$ws.Range("N80").Value = -38168.39999999999
$ws.Range("H83").Value = 12168.454
$ws.Range("J83").Value = 12098.8
$ws.Range("L83").Value = 108889.2
$ws.Range("N83").Value = -118249.2
$ws.Range("H111").Value = 106.666664
$ws.Range("I111").Value = 106.666664
$ws.Range("K111").Value = 319.999992
$ws.Range("M111").Value = 2747.000008
$ws.Range("H119").Value = 219
$ws.Range("I119").Value = 219
$ws.Range("K119").Value = 657
$ws.Range("M119").Value = 4181
$ws.Range("H131").Value = 1093
$ws.Range("J131").Value = 1499
$ws.Range("L131").Value = 4497
$ws.Range("N131").Value = -14577

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2130.3333
$ws.Range("I132").Value = 2153.6
$ws.Range("K132").Value = 6460.799999999999
$ws.Range("M132").Value = -3930.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 819.4666999999999
$ws.Range("I22").Value = 885.63635
$ws.Range("J22").Value = 637.5
$ws.Range("K22").Value = 885.63635
$ws.Range("L22").Value = 637.5
$ws.Range("M22").Value = -590.63635
$ws.Range("N22").Value = -1227.5
$ws.Range("H27").Value = 819.4666999999999
$ws.Range("I27").Value = 885.63635
$ws.Range("J27").Value = 637.5
$ws.Range("K27").Value = 885.63635
$ws.Range("L27").Value = 637.5
$ws.Range("M27").Value = -778.63635
$ws.Range("N27").Value = -851.5
$ws.Range("H55").Value = 234.86667
$ws.Range("I55").Value = 154.22223
$ws.Range("J55").Value = 355.83334
$ws.Range("K55").Value = 154.22223
$ws.Range("L55").Value = 355.83334
$ws.Range("M55").Value = 18.77777
$ws.Range("N55").Value = -701.83334
$ws.Range("H132").Value = 6255.222
$ws.Range("I132").Value = 5700.875
$ws.Range("J132").Value = 6698.7
$ws.Range("K132").Value = 17102.625
$ws.Range("L132").Value = 20096.1
$ws.Range("M132").Value = -14572.625
$ws.Range("N132").Value = -25156.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6666.3335
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H68").Value = 54950
$ws.Range("J68").Value = 54950
$ws.Range("L68").Value = 54950
$ws.Range("N68").Value = -56572
$ws.Range("H71").Value = 54950
$ws.Range("J71").Value = 54950
$ws.Range("L71").Value = 164850
$ws.Range("N71").Value = -172962
$ws.Range("H132").Value = 944.6667
$ws.Range("I132").Value = 882.125
$ws.Range("J132").Value = 1016.1429
$ws.Range("K132").Value = 2646.375
$ws.Range("M132").Value = -116.375
$ws.Range("N132").Value = -8108.4287
